# Make the "Token" column header cells bold for the CONS, CONSEN, IDEN,
# SALIDA, FRASE and ENTRADA rows in the Tokens table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$tokens = @("CONS", "CONSEN", "IDEN", "SALIDA", "FRASE", "ENTRADA")

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Cell($i, 1)
    $text = $cell.Range.Text
    # Cell range text includes trailing cell-mark characters; trim them.
    $text = $text.TrimEnd([char]7, [char]13)
    if ($tokens -contains $text) {
        $cell.Range.Font.Bold = $true
    }
}
